$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("P2").Value = 2.1
$ws.Range("T2").Value = 1.86
$ws.Range("X2").Value = 16.5

# Row 4 updates
$ws.Range("F4").Value = 8.6
$ws.Range("G4").Value = 10.5
$ws.Range("H4").Value = 1.42
$ws.Range("I4").Value = 1.45
$ws.Range("Q4").Value = 1.88
$ws.Range("S4").Value = 3.2
$ws.Range("V4").Value = 3.2
$ws.Range("W4").Value = 1.11
$ws.Range("X4").Value = 23
$ws.Range("Z4").Value = 8
$ws.Range("AB4").Value = 28
$ws.Range("AF4").Value = 90
$ws.Range("AJ4").Value = 410
$ws.Range("AK4").Value = 190
$ws.Range("AL4").Value = 160
$ws.Range("AN4").Value = 300
